$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kosten")

# Updated / new data rows (Zeit, Beste Werte, Optimale Werte)
$data = @(
    @(5,  28600, 3400),
    @(10, 25000, 3400),
    @(15, 24300, 3400),
    @(20, 24000, 3400),
    @(25, 24000, 3400),
    @(30, 24000, 3400),
    @(35, 24000, 3400),
    @(40, 24000, 3400),
    @(45, 23800, 3400),
    @(50, 23800, 3400),
    @(55, 23800, 3400),
    @(60, 23800, 3400),
    @(65, 23800, 3400)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}

# Update chart series to cover the expanded data range
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$chart.SeriesCollection().Item(1).Formula = "=SERIES(Kosten!`$B`$1,Kosten!`$A`$2:`$A`$14,Kosten!`$B`$2:`$B`$14,1)"
$chart.SeriesCollection().Item(2).Formula = "=SERIES(Kosten!`$C`$1,Kosten!`$A`$2:`$A`$14,Kosten!`$C`$2:`$C`$14,2)"
